# CORE_holdings.xlsx - "Add files via upload" update
# - Bumps the "as of" date in the confidential disclosure text (A11) from
#   2021-05-11 to 2021-05-12.
# - Refreshes the Weight (D) and Percent Change (E) figures for rows 2-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; temporarily unprotect so the cells can be
# written, then restore protection afterwards.
$ws.Unprotect()

$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-12 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.5006285541645973
$ws.Range("E2").Value = -0.01736555159184217

$ws.Range("D3").Value = 0.2434991646720012
$ws.Range("E3").Value = -0.02498152254249819

$ws.Range("D4").Value = 0.09567425309653389
$ws.Range("E4").Value = -0.02861825516893601

$ws.Range("D5").Value = 0.1036549572686555
$ws.Range("E5").Value = -0.03084595493165865

$ws.Range("D6").Value = 0.02980702257549064
$ws.Range("E6").Value = -0.0310241834473457

$ws.Range("D7").Value = 0.02673604822272146
$ws.Range("E7").Value = -0.03339897597479324

$ws.Range("E8").Value = -0.02252973235424649

$ws.Protect($null, $true, $true, $true)
